$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 21, shifting rows 21-23 down to 22-24
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with the $GPVTG NMEA sentence header and two "x" placeholder cells
$ws.Range("A21").Value = '$GPVTG'
$ws.Range("L21").Value = "x"
$ws.Range("M21").Value = "x"

# Clear the blank cells that copied formatting from the row above during the insert
$ws.Range("B21:K21").Clear()

# Update the active selection to reflect where the user left off editing
$ws.Range("P27").Select()
